# Generate Report for Handoff
# Replaces the source-file UUID / handoff artifact names and timestamps with the
# values produced by the newer handoff run, and clears the (now stale) "Latest
# Target File" / "Latest Handback File" columns on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "84191805-2b27-4f11-8519-64ea8ae9ae5f"
$newGuid = "37418fde-e613-466b-b38b-638ec7b63a5a"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"

$oldHoHash = "5820819b002e23390a75775157cee49df6858b01"
$newHoHash = "577f7a8308fe14cb74e729602c36d10eb8b88b8f"

# Hyperlink target addresses are unchanged by this commit - only the visible
# text that mirrors the source file name is refreshed.
$sourceHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbfc818b29d5655fa75dcdef76cbf44c62203415/e2e/$oldMdName"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newMdName

$ws.Range("B1").Hyperlinks.Delete()
$ws.Range("B2").Value = "e2e\$newMdName"
$ws.Hyperlinks.Add($ws.Range("B2"), $sourceHyperlinkUrl, "", "", "e2e\$newMdName")
$ws.Range("B2").Style = "Hyperlink"

$ws.Range("G2").Value = "2016-08-26 20:58:22"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()
$ws.Range("A2").Value = $newMdName
$ws.Hyperlinks.Add($ws.Range("A2"), $sourceHyperlinkUrl, "", "", $newMdName)
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("G2").Value = "$newGuid.$newHoHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-26 20:58:17"

# Latest Target File / Latest Handback File are no longer available for this
# handoff - clear the values and drop the now-removed hyperlink on I2.
$ws.Range("I2").ClearContents()
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""

$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(9).ColumnWidth = 17.8
$ws.Columns.Item(10).ColumnWidth = 20.9

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()
$ws.Range("A2").Value = $newMdName
$ws.Hyperlinks.Add($ws.Range("A2"), $sourceHyperlinkUrl, "", "", $newMdName)
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("G2").Value = "$newGuid.$newHoHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-26 20:58:22"

$ws.Range("I2").ClearContents()
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""

$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(9).ColumnWidth = 17.8
$ws.Columns.Item(10).ColumnWidth = 20.9

Write-Output "Localization status report refreshed for $newGuid"
